$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet title and update running total label text (now "through December 10")
$wb.Worksheets.Item(1).Name = "Through 2022-12-10"
$ws.Range("B1").Value = "December 2022 (through December 10)"

# Update/insert cell values per new data for 2022-12-10
$ws.Range("AX2").Value = 3

$ws.Range("B3").Value = 1
$ws.Range("Z3").Value = 1

$ws.Range("N4").Value = 4
$ws.Range("Z4").Value = 2

$ws.Range("AL12").Value = 1

$ws.Range("B14").Value = 3
$ws.Range("Z14").Value = 5

$ws.Range("B15").Value = 4

$ws.Range("B16").Value = 1

$ws.Range("B18").Value = 2
$ws.Range("BV18").Value = 1

$ws.Range("N20").Value = 3

$ws.Range("Z23").Value = 3

$ws.Range("Z28").Value = 1

$ws.Range("B29").Value = 1
$ws.Range("BV29").Value = 1

$ws.Range("CH36").Value = 1

$ws.Range("AX40").Value = 2

$ws.Range("AX51").Value = 1

$ws.Range("N67").Value = 1

$ws.Range("N75").Value = 1
